$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing key/value pairs currently in rows 331-342 (z3001..z3012)
# down to rows 341-352, to make room for the new z2009 pair at row 329.
# Copy bottom-up so we don't overwrite source rows before they've been read.
for ($i = 12; $i -ge 1; $i--) {
    $srcRow = 330 + $i
    $dstRow = 340 + $i
    $ws.Cells.Item($dstRow, 1).Value2 = $ws.Cells.Item($srcRow, 1).Value2
    $ws.Cells.Item($dstRow, 2).Value2 = $ws.Cells.Item($srcRow, 2).Value2
}

# Remove the now-duplicated old rows 331-340 entirely (no leftover formatting)
$ws.Range("A331:B340").Clear()

# Add the new key/value pair for z2009 at row 329
$ws.Cells.Item(329, 1).Value2 = "z2009"
$ws.Cells.Item(329, 2).Value2 = "평면 벡터들 사이의 관계식으로 부터 점들의 위치를 나타내는 도형을 파악해서 문제에서 요구하는 최대값, 최솟값을 구합니다. "

# Restore the user's selection reflected in the saved file
$ws.Range("B332").Select() | Out-Null
